# Applies the "Add files via upload" data refresh to the VIC key outbreaks sheet:
# replaces the Cluster name / Active cases / Exist rows (A2:C44) with the refreshed
# 59-row dataset (A2:C60), matching the new authoritative extract.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is (Cluster name, Active cases, Exist)
$data = @(
  ,("139 Highett St Apartment Complex Richmond", 10, "old")
  ,("3175 The Bays Aged Care Facility Hastings", 16, "new")
  ,("3175 The Bays Aged Care Facility Hastings", 17, "old")
  ,("3563 Embracia Aged Care Reservoir", 11, "new")
  ,("3600 Belvedere Age Care Noble Park Outbreak", 21, "old")
  ,("3600 Belvedere Aged Care Noble Park", 22, "new")
  ,("3612 BlueCross Glengowrie", 35, "new")
  ,("3612 BlueCross Glengowrie Outbreak", 33, "old")
  ,("3684 Homestyle Aged Care Langford Grange Cranbourne East", 25, "new")
  ,("3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak", 23, "old")
  ,("3980 Arcare Keysborough Aged Care Keysborough", 13, "new")
  ,("4075 Ferndale Gardens Aged Care Services Bayswater North", 17, "new")
  ,("4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak", 17, "old")
  ,("4518 Regis Aged Care Fawkner", 14, "new")
  ,("4518 Regis Aged Care Fawkner Outbreak", 13, "old")
  ,("AW Window Transport Group Depot North Geelong", 10, "new")
  ,("AW Window Transport Group Depot North Geelong Outbreak", 10, "old")
  ,("Allied Pinnacle Factory Altona North", 12, "new")
  ,("Bread Solutions Braeside", 15, "new")
  ,("Bread Solutions Braeside Outbreak", 18, "old")
  ,("CS Square Caroline Springs", 13, "new")
  ,("CS Square Caroline Springs Outbreak", 14, "old")
  ,("Cedar Meats Australia Brooklyn", 10, "new")
  ,("Cedar Meats Australia Brooklyn Outbreak", 10, "old")
  ,("Comfort Sleep Bedding Co Thomastown", 10, "new")
  ,("Community Kids Pascoe Vale Early Education Centre Pascoe Vale", 22, "new")
  ,("Community Kids Pascoe Vale Early Education Centre Pascoe Vale Outbreak", 20, "old")
  ,("Embracia Aged Care Reservoir Outbreak", 19, "old")
  ,("Guardian Childcare Caulfield", 20, "new")
  ,("Guardian Childcare Caulfield Outbreak", 20, "old")
  ,("Hello Fresh Warehouse Ravenhall", 13, "old")
  ,("Hello Fresh Warehouse Ravenhall", 17, "new")
  ,("Inghams Enterprise Somerville Outbreak", 17, "old")
  ,("Inghams Enterprises Somerville", 24, "new")
  ,("Kool Kidz Childcare Narre Warren", 10, "old")
  ,("Lantmannen Unibake Australia Mordialloc", 22, "new")
  ,("Lantmannen Unibake Australia Mordialloc", 25, "old")
  ,("Monash Health Kingston Centre South 5", 10, "new")
  ,("Nido Early School Ascot Vale", 11, "old")
  ,("Northern Health The Northern Hospital Epping", 14, "old")
  ,("Northern Health The Northern Hospital Epping", 17, "new")
  ,("Oceania Meat Processors Laverton North", 16, "new")
  ,("Oceania Meat Processors Laverton North Outbreak", 16, "old")
  ,("Pick It Up Fitness Mulgrave Outbreak", 11, "old")
  ,("Robin Hood Inn Drouin West Outbreak", 42, "old")
  ,("St Vincents Hospital Emergency Department Melbourne", 36, "new")
  ,("St Vincents Hospital Emergency Department Melbourne", 41, "old")
  ,("Target Distribution Centre Truganina", 17, "new")
  ,("Target Distribution Centre Truganina Outbreak", 19, "old")
  ,("The Robin Hood Inn Drouin West", 46, "new")
  ,("The Royal Children`'s Hospital Melbourne Emergency Department Parkville Tier 1B", 12, "old")
  ,("The Royal Melbourne Hospital Parkville Emergency Department", 10, "new")
  ,("The Toolshed Bar Private Event Noojee", 16, "new")
  ,("The Toolshed Bar Private Event Noojee", 17, "old")
  ,("Werribee Mercy Hospital Emergency Department", 23, "new")
  ,("Werribee Mercy Hospital Emergency Department", 25, "old")
  ,("Western Health Footscray Hospital Emergency Department", 10, "new")
  ,("Western Health Sunshine Hospital Emergency Department", 21, "old")
  ,("Western Health Sunshine Hospital Emergency Department", 22, "new")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $item = $data[$i]
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $item[0]
  $ws.Cells.Item($r, 2).Value = $item[1]
  $ws.Cells.Item($r, 3).Value = $item[2]
}
